$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = 2.3
$ws.Range("H5").Value = 3.5
$ws.Range("I5").Value = 2.88
$ws.Range("J5").Value = 2.88
$ws.Range("L5").Value = 3.4
$ws.Range("M5").Value = 1.04
$ws.Range("N5").Value = 13
$ws.Range("X5").Value = 13
$ws.Range("AC5").Value = 13
$ws.Range("AK5").Value = 29
$ws.Range("AS5").Value = 126
